$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "69.651.05"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "3.641.06"
$ws.Range("E3").Value = "  +18.35%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +2.64%  "
Set-TextValue "D6" "186.13"
$ws.Range("E6").Value = "  +8.77%  "
$ws.Range("D7").Value = "3.640.52"
$ws.Range("E7").Value = "  +18.43%  "
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +4.59%  "
Set-TextValue "D10" "0.162"
$ws.Range("E10").Value = "  +8.80%  "
$ws.Range("E11").Value = "  +4.87%  "
Set-TextValue "D12" "0.497"
$ws.Range("E12").Value = "  +5.79%  "
Set-TextValue "D13" "39.34"
$ws.Range("E13").Value = "  +9.98%  "
Set-TextValue "D14" "0.0000253"
$ws.Range("E14").Value = "  +5.93%  "
$ws.Range("D15").Value = "4.260.80"
$ws.Range("E15").Value = "  +18.63%  "
$ws.Range("D16").Value = "3.650.03"
$ws.Range("E16").Value = "  +18.58%  "
$ws.Range("D17").Value = "69.880.88"
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("E18").Value = "  +1.96%  "
Set-TextValue "D19" "7.52"
$ws.Range("E19").Value = "  +8.13%  "
Set-TextValue "D20" "17.19"
$ws.Range("E20").Value = "  +3.66%  "
Set-TextValue "D21" "508.90"
$ws.Range("E21").Value = "  +4.56%  "
Set-TextValue "D22" "9.31"
$ws.Range("E22").Value = "  +21.15%  "
Set-TextValue "D23" "0.749"
$ws.Range("E23").Value = "  +9.16%  "
Set-TextValue "D24" "87.90"
$ws.Range("E24").Value = "  +6.65%  "
Set-TextValue "D25" "13.56"
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("E26").Value = "  +8.42%  "
Set-TextValue "D27" "10.89"
$ws.Range("E27").Value = "  +7.68%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +13.09%  "
Set-TextValue "D30" "8.20"
$ws.Range("E30").Value = "  +4.55%  "
Set-TextValue "D31" "32.33"
$ws.Range("E31").Value = "  +16.62%  "
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("E33").Value = "  +18.73%  "
$ws.Range("E34").Value = "  +5.39%  "
$ws.Range("E35").Value = "  +0.06%  "
Set-TextValue "D36" "6.18"
$ws.Range("E36").Value = "  +11.11%  "
Set-TextValue "D37" "1.02"
$ws.Range("E37").Value = "  +8.44%  "
$ws.Range("E38").Value = "  +11.09%  "
Set-TextValue "D39" "2.10"
$ws.Range("E39").Value = "  +7.17%  "
Set-TextValue "D40" "46.85"
$ws.Range("E40").Value = "  -2.36%  "
Set-TextValue "D41" "50.71"
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("D43").Value = "3.164.17"
$ws.Range("E43").Value = "  +14.16%  "
$ws.Range("E44").Value = "  +7.30%  "
Set-TextValue "D45" "2.79"
$ws.Range("E45").Value = "  +10.41%  "
Set-TextValue "D46" "405.54"
$ws.Range("E46").Value = "  +11.06%  "
$ws.Range("E47").Value = "  +6.77%  "
Set-TextValue "D48" "27.84"
$ws.Range("E48").Value = "  +14.64%  "
Set-TextValue "D49" "136.74"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D50" "1.00"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D51" "2.45"
$ws.Range("E51").Value = "  +13.86%  "
